# Weekly price-list update: a new weekly observation is inserted at row 95
# (the Excel row immediately after the header), pushing every existing
# record down by one row (95->96, 96->97, ... 231->232).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 95; Excel shifts rows 95..231 down to 96..232
# and the sheet's used range grows from A1:R231 to A1:R232 automatically.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(95, 1).Value  = 8
$ws.Cells.Item(95, 2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(95, 3).Value  = 'Coquimbo'
$ws.Cells.Item(95, 4).Value  = 44799
$ws.Cells.Item(95, 5).Value  = 4
$ws.Cells.Item(95, 6).Value  = 100112037
$ws.Cells.Item(95, 7).Value  = 'Cebollín'
$ws.Cells.Item(95, 8).Value  = 'Sin especificar'
$ws.Cells.Item(95, 9).Value  = 'Primera'
$ws.Cells.Item(95, 10).Value = 2000
$ws.Cells.Item(95, 11).Value = 1400
$ws.Cells.Item(95, 12).Value = 1600
$ws.Cells.Item(95, 13).Value = 1500
$ws.Cells.Item(95, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(95, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(95, 16).Value = 250
$ws.Cells.Item(95, 17).Value = 6
$ws.Cells.Item(95, 18).Value = 'Hortaliza'
